$d = $word.ActiveDocument

$replacements = @(
    @{old="87×18="; new="38×78="},
    @{old="70×60="; new="68×42="},
    @{old="78×73="; new="13×46="},
    @{old="27×51="; new="34×85="},
    @{old="38×59="; new="89×98="},
    @{old="57×77="; new="28×68="},
    @{old="71×16="; new="66×12="},
    @{old="59×60="; new="58×51="},
    @{old="95×53="; new="55×49="},
    @{old="45×81="; new="87×99="},
    @{old="31×43="; new="75×21="},
    @{old="39×78="; new="98×50="},
    @{old="30×55="; new="19×42="},
    @{old="58×47="; new="85×63="},
    @{old="95×71="; new="68×35="},
    @{old="57×98="; new="23×34="},
    @{old="13×44="; new="48×87="},
    @{old="16×88="; new="23×80="},
    @{old="33×14="; new="25×45="},
    @{old="87×41="; new="64×19="},
    @{old="25×91="; new="29×38="},
    @{old="46×46="; new="58×37="},
    @{old="98×32="; new="62×31="},
    @{old="27×69="; new="66×64="},
    @{old="32×40="; new="14×72="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
